$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-25 21:48:30'
$ws.Range('O2').Value = '5.4 °C'
$ws.Range('E3').Value = '2026-02-25 21:48:33'
$ws.Range('H3').NumberFormat = '@'
$ws.Range('H3').Value = '38%'
$ws.Range('N3').Value = '1.1 °C 21:20 TU'
$ws.Range('O3').Value = '3.7 °C'
$ws.Range('E4').Value = '2026-02-25 21:48:35'
$ws.Range('J4').Value = '1022.0 hPa'
$ws.Range('O4').Value = '8.8 °C'
$ws.Range('E5').Value = '2026-02-25 21:48:37'
$ws.Range('H5').NumberFormat = '@'
$ws.Range('H5').Value = '30%'
$ws.Range('N5').Value = '1.6 °C 21:29 TU'
$ws.Range('O5').Value = '5.6 °C'
$ws.Range('E6').Value = '2026-02-25 21:48:40'
$ws.Range('J6').Value = '1021.9 hPa'
$ws.Range('E7').Value = '2026-02-25 21:48:42'
$ws.Range('J7').Value = '1021.5 hPa'
$ws.Range('E8').Value = '2026-02-25 21:48:45'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '80%'
$ws.Range('J8').Value = '1021.2 hPa'
$ws.Range('O8').Value = '11.8 °C'
$ws.Range('E9').Value = '2026-02-25 21:48:47'
$ws.Range('O9').Value = '10.3 °C'
$ws.Range('E10').Value = '2026-02-25 21:48:50'
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value = '90%'
$ws.Range('O10').Value = '9.3 °C'
$ws.Range('E11').Value = '2026-02-25 21:48:52'
$ws.Range('O11').Value = '8.8 °C'
$ws.Range('E12').Value = '2026-02-25 21:48:55'
$ws.Range('O12').Value = '10.3 °C'
$ws.Range('E13').Value = '2026-02-25 21:48:57'
$ws.Range('J13').Value = '1022.9 hPa'
$ws.Range('O13').Value = '6.7 °C'
$ws.Range('E14').Value = '2026-02-25 21:49:00'
$ws.Range('O14').Value = '10.7 °C'
$ws.Range('E15').Value = '2026-02-25 21:49:02'
$ws.Range('E16').Value = '2026-02-25 21:49:04'
$ws.Range('O16').Value = '3.1 °C'
$ws.Range('E17').Value = '2026-02-25 21:49:07'
$ws.Range('O17').Value = '9.1 °C'
$ws.Range('E18').Value = '2026-02-25 21:49:09'
$ws.Range('J18').Value = '1022.1 hPa'
$ws.Range('E19').Value = '2026-02-25 21:49:12'
$ws.Range('H19').NumberFormat = '@'
$ws.Range('H19').Value = '50%'
$ws.Range('N19').Value = '7.7 °C 21:12 TU'
$ws.Range('O19').Value = '12.1 °C'
$ws.Range('E20').Value = '2026-02-25 21:49:14'
$ws.Range('N20').Value = '-1.9 °C 21:24 TU'
$ws.Range('O20').Value = '2.7 °C'
$ws.Range('E21').Value = '2026-02-25 21:49:17'
$ws.Range('J21').Value = '1021.6 hPa'
$ws.Range('E22').Value = '2026-02-25 21:49:19'
$ws.Range('N22').Value = '0.1 °C 21:28 TU'
$ws.Range('O22').Value = '2.4 °C'
$ws.Range('E23').Value = '2026-02-25 21:49:22'
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = '35%'
$ws.Range('N23').Value = '1.5 °C 21:08 TU'
$ws.Range('O23').Value = '3.9 °C'
$ws.Range('E24').Value = '2026-02-25 21:49:24'
$ws.Range('J24').Value = '1020.2 hPa'
$ws.Range('L24').Value = '17.6 km/h - 90º 21:29 TU'
$ws.Range('O24').Value = '10.9 °C'
$ws.Range('E25').Value = '2026-02-25 21:49:27'
$ws.Range('E26').Value = '2026-02-25 21:49:29'
$ws.Range('H26').NumberFormat = '@'
$ws.Range('H26').Value = '48%'
$ws.Range('J26').Value = '1019.8 hPa'
$ws.Range('N26').Value = '5.2 °C 21:09 TU'
$ws.Range('O26').Value = '9.9 °C'
$ws.Range('E27').Value = '2026-02-25 21:49:32'
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = '42%'
$ws.Range('N27').Value = '1.9 °C 21:29 TU'
$ws.Range('O27').Value = '5.2 °C'
$ws.Range('E28').Value = '2026-02-25 21:49:34'
$ws.Range('E29').Value = '2026-02-25 21:49:37'
$ws.Range('O29').Value = '11.9 °C'
$ws.Range('E30').Value = '2026-02-25 21:49:39'
$ws.Range('J30').Value = '1022.0 hPa'
$ws.Range('E31').Value = '2026-02-25 21:49:42'
$ws.Range('J31').Value = '1021.6 hPa'
$ws.Range('E32').Value = '2026-02-25 21:49:44'
$ws.Range('O32').Value = '9.1 °C'
$ws.Range('E33').Value = '2026-02-25 21:49:47'
$ws.Range('E34').Value = '2026-02-25 21:49:49'
$ws.Range('O34').Value = '3.3 °C'
$ws.Range('E35').Value = '2026-02-25 21:49:52'
$ws.Range('J35').Value = '1019.7 hPa'
$ws.Range('O35').Value = '12.4 °C'
$ws.Range('E36').Value = '2026-02-25 21:49:54'
$ws.Range('J36').Value = '1022.1 hPa'
$ws.Range('E37').Value = '2026-02-25 21:49:57'
$ws.Range('J37').Value = '1023.6 hPa'
$ws.Range('E38').Value = '2026-02-25 21:49:59'
$ws.Range('E39').Value = '2026-02-25 21:50:02'
$ws.Range('E40').Value = '2026-02-25 21:50:04'
$ws.Range('H40').NumberFormat = '@'
$ws.Range('H40').Value = '61%'
$ws.Range('J40').Value = '1022.0 hPa'
$ws.Range('O40').Value = '9.5 °C'
$ws.Range('E41').Value = '2026-02-25 21:50:06'
$ws.Range('J41').Value = '1021.1 hPa'
$ws.Range('E42').Value = '2026-02-25 21:50:09'
$ws.Range('E43').Value = '2026-02-25 21:50:11'
$ws.Range('O43').Value = '9.8 °C'
$ws.Range('E44').Value = '2026-02-25 21:50:14'
$ws.Range('H44').NumberFormat = '@'
$ws.Range('H44').Value = '46%'
$ws.Range('E45').Value = '2026-02-25 21:50:16'
$ws.Range('J45').Value = '1020.0 hPa'
$ws.Range('O45').Value = '10.8 °C'
$ws.Range('E46').Value = '2026-02-25 21:50:19'
$ws.Range('J46').Value = '1020.9 hPa'
